$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 163, shifting existing rows 163:170 down to 164:171
$ws.Rows.Item(163).Insert()

# Populate the newly inserted row 163 with the new weekly record
$ws.Cells.Item(163, 1).Value = 6
$ws.Cells.Item(163, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(163, 3).Value = "Metropolitana"
$ws.Cells.Item(163, 4).Value = 44610
$ws.Cells.Item(163, 5).Value = 13
$ws.Cells.Item(163, 6).Value = "Fruta"
$ws.Cells.Item(163, 7).Value = 100101
$ws.Cells.Item(163, 8).Value = "Berries"
$ws.Cells.Item(163, 9).Value = 100101004
$ws.Cells.Item(163, 10).Value = "Frambuesa"
$ws.Cells.Item(163, 11).Value = "Sin especificar"
$ws.Cells.Item(163, 12).Value = "Especial"
$ws.Cells.Item(163, 13).Value = 350
$ws.Cells.Item(163, 14).Value = 8000
$ws.Cells.Item(163, 15).Value = 8000
$ws.Cells.Item(163, 16).Value = 8000
$ws.Cells.Item(163, 17).Value = '$/bandeja 2 kilos'
$ws.Cells.Item(163, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(163, 19).Value = 4000
$ws.Cells.Item(163, 20).Value = 2
